$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing coefficient values in column B
$ws.Range("B2").Value = -0.6161452981289695
$ws.Range("B3").Value = 0.9412207439789471
$ws.Range("B4").Value = 588.5888204728865

# Delete row 5 entirely (removes A5 "4" label and its B5 value)
$ws.Rows.Item(5).Delete()
